$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-08-19 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-20 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("93-48=45", $true, $false, $false, $false, $false, $true, 1, $false, "46-44=2", 2) | Out-Null
$d.Content.Find.Execute("5+38=43", $true, $false, $false, $false, $false, $true, 1, $false, "63-58=5", 2) | Out-Null
$d.Content.Find.Execute("35-21=14", $true, $false, $false, $false, $false, $true, 1, $false, "47+25=72", 2) | Out-Null
$d.Content.Find.Execute("83-62=21", $true, $false, $false, $false, $false, $true, 1, $false, "69+29=98", 2) | Out-Null
$d.Content.Find.Execute("0+39=39", $true, $false, $false, $false, $false, $true, 1, $false, "61-51=10", 2) | Out-Null
$d.Content.Find.Execute("42+44=86", $true, $false, $false, $false, $false, $true, 1, $false, "96-30=66", 2) | Out-Null
$d.Content.Find.Execute("97-62=35", $true, $false, $false, $false, $false, $true, 1, $false, "58-23=35", 2) | Out-Null
$d.Content.Find.Execute("26-10=16", $true, $false, $false, $false, $false, $true, 1, $false, "45-42=3", 2) | Out-Null
$d.Content.Find.Execute("87-42=45", $true, $false, $false, $false, $false, $true, 1, $false, "81-9=72", 2) | Out-Null
$d.Content.Find.Execute("58-21=37", $true, $false, $false, $false, $false, $true, 1, $false, "70-16=54", 2) | Out-Null
$d.Content.Find.Execute("18-2=16", $true, $false, $false, $false, $false, $true, 1, $false, "38-13=25", 2) | Out-Null
$d.Content.Find.Execute("61-36=25", $true, $false, $false, $false, $false, $true, 1, $false, "57-18=39", 2) | Out-Null
$d.Content.Find.Execute("22+30=52", $true, $false, $false, $false, $false, $true, 1, $false, "66-47=19", 2) | Out-Null
$d.Content.Find.Execute("31+4=35", $true, $false, $false, $false, $false, $true, 1, $false, "6+5=11", 2) | Out-Null
$d.Content.Find.Execute("60+38=98", $true, $false, $false, $false, $false, $true, 1, $false, "57-15=42", 2) | Out-Null
$d.Content.Find.Execute("11+41=52", $true, $false, $false, $false, $false, $true, 1, $false, "40+0=40", 2) | Out-Null
$d.Content.Find.Execute("91-46=45", $true, $false, $false, $false, $false, $true, 1, $false, "62-34=28", 2) | Out-Null
$d.Content.Find.Execute("12-6=6", $true, $false, $false, $false, $false, $true, 1, $false, "37-29=8", 2) | Out-Null
$d.Content.Find.Execute("64+9=73", $true, $false, $false, $false, $false, $true, 1, $false, "63+22=85", 2) | Out-Null
$d.Content.Find.Execute("68-13=55", $true, $false, $false, $false, $false, $true, 1, $false, "23+68=91", 2) | Out-Null
$d.Content.Find.Execute("6+41=47", $true, $false, $false, $false, $false, $true, 1, $false, "17+19=36", 2) | Out-Null
$d.Content.Find.Execute("69+3=72", $true, $false, $false, $false, $false, $true, 1, $false, "81-6=75", 2) | Out-Null
$d.Content.Find.Execute("85-57=28", $true, $false, $false, $false, $false, $true, 1, $false, "81-28=53", 2) | Out-Null
$d.Content.Find.Execute("38-1=37", $true, $false, $false, $false, $false, $true, 1, $false, "60+0=60", 2) | Out-Null
$d.Content.Find.Execute("60-19=41", $true, $false, $false, $false, $false, $true, 1, $false, "10+27=37", 2) | Out-Null
$d.Content.Find.Execute("76-64=12", $true, $false, $false, $false, $false, $true, 1, $false, "87+7=94", 2) | Out-Null
$d.Content.Find.Execute("30+19=49", $true, $false, $false, $false, $false, $true, 1, $false, "88-31=57", 2) | Out-Null
$d.Content.Find.Execute("40+3=43", $true, $false, $false, $false, $false, $true, 1, $false, "69-23=46", 2) | Out-Null
$d.Content.Find.Execute("23-7=16", $true, $false, $false, $false, $false, $true, 1, $false, "11+6=17", 2) | Out-Null
$d.Content.Find.Execute("75+21=96", $true, $false, $false, $false, $false, $true, 1, $false, "3+3=6", 2) | Out-Null
$d.Content.Find.Execute("87-55=32", $true, $false, $false, $false, $false, $true, 1, $false, "98-4=94", 2) | Out-Null
$d.Content.Find.Execute("9+85=94", $true, $false, $false, $false, $false, $true, 1, $false, "47-2=45", 2) | Out-Null
$d.Content.Find.Execute("42+22=64", $true, $false, $false, $false, $false, $true, 1, $false, "54+19=73", 2) | Out-Null
$d.Content.Find.Execute("97-74=23", $true, $false, $false, $false, $false, $true, 1, $false, "3-0=3", 2) | Out-Null
$d.Content.Find.Execute("59+16=75", $true, $false, $false, $false, $false, $true, 1, $false, "69-2=67", 2) | Out-Null
$d.Content.Find.Execute("29+3=32", $true, $false, $false, $false, $false, $true, 1, $false, "8+1=9", 2) | Out-Null
$d.Content.Find.Execute("94-50=44", $true, $false, $false, $false, $false, $true, 1, $false, "31-4=27", 2) | Out-Null
$d.Content.Find.Execute("56-24=32", $true, $false, $false, $false, $false, $true, 1, $false, "43+10=53", 2) | Out-Null
$d.Content.Find.Execute("25-1=24", $true, $false, $false, $false, $false, $true, 1, $false, "86-76=10", 2) | Out-Null
$d.Content.Find.Execute("74-51=23", $true, $false, $false, $false, $false, $true, 1, $false, "45-2=43", 2) | Out-Null
$d.Content.Find.Execute("2+35=37", $true, $false, $false, $false, $false, $true, 1, $false, "65-7=58", 2) | Out-Null
$d.Content.Find.Execute("24+17=41", $true, $false, $false, $false, $false, $true, 1, $false, "67+11=78", 2) | Out-Null
$d.Content.Find.Execute("23+62=85", $true, $false, $false, $false, $false, $true, 1, $false, "78-77=1", 2) | Out-Null
$d.Content.Find.Execute("33+7=40", $true, $false, $false, $false, $false, $true, 1, $false, "80-51=29", 2) | Out-Null
$d.Content.Find.Execute("74+16=90", $true, $false, $false, $false, $false, $true, 1, $false, "4+91=95", 2) | Out-Null
$d.Content.Find.Execute("28+59=87", $true, $false, $false, $false, $false, $true, 1, $false, "32-16=16", 2) | Out-Null
$d.Content.Find.Execute("89-30=59", $true, $false, $false, $false, $false, $true, 1, $false, "57-42=15", 2) | Out-Null
$d.Content.Find.Execute("67+0=67", $true, $false, $false, $false, $false, $true, 1, $false, "1+34=35", 2) | Out-Null
$d.Content.Find.Execute("21+7=28", $true, $false, $false, $false, $false, $true, 1, $false, "60-29=31", 2) | Out-Null
$d.Content.Find.Execute("35-34=1", $true, $false, $false, $false, $false, $true, 1, $false, "26+43=69", 2) | Out-Null
$d.Content.Find.Execute("16+80=96", $true, $false, $false, $false, $false, $true, 1, $false, "29+40=69", 2) | Out-Null
$d.Content.Find.Execute("5+53=58", $true, $false, $false, $false, $false, $true, 1, $false, "52+39=91", 2) | Out-Null
$d.Content.Find.Execute("25-9=16", $true, $false, $false, $false, $false, $true, 1, $false, "75-62=13", 2) | Out-Null
$d.Content.Find.Execute("57+38=95", $true, $false, $false, $false, $false, $true, 1, $false, "93-35=58", 2) | Out-Null
$d.Content.Find.Execute("47+31=78", $true, $false, $false, $false, $false, $true, 1, $false, "40+46=86", 2) | Out-Null
$d.Content.Find.Execute("83-54=29", $true, $false, $false, $false, $false, $true, 1, $false, "4+39=43", 2) | Out-Null
$d.Content.Find.Execute("8+55=63", $true, $false, $false, $false, $false, $true, 1, $false, "31-1=30", 2) | Out-Null
$d.Content.Find.Execute("44-27=17", $true, $false, $false, $false, $false, $true, 1, $false, "80-25=55", 2) | Out-Null
$d.Content.Find.Execute("27+72=99", $true, $false, $false, $false, $false, $true, 1, $false, "92+0=92", 2) | Out-Null
$d.Content.Find.Execute("77-4=73", $true, $false, $false, $false, $false, $true, 1, $false, "1+50=51", 2) | Out-Null
$d.Content.Find.Execute("71-11=60", $true, $false, $false, $false, $false, $true, 1, $false, "33-31=2", 2) | Out-Null
$d.Content.Find.Execute("81+12=93", $true, $false, $false, $false, $false, $true, 1, $false, "78-52=26", 2) | Out-Null
$d.Content.Find.Execute("88-38=50", $true, $false, $false, $false, $false, $true, 1, $false, "54+14=68", 2) | Out-Null
$d.Content.Find.Execute("46-35=11", $true, $false, $false, $false, $false, $true, 1, $false, "39+54=93", 2) | Out-Null
$d.Content.Find.Execute("4+22=26", $true, $false, $false, $false, $false, $true, 1, $false, "80+5=85", 2) | Out-Null
$d.Content.Find.Execute("4+54=58", $true, $false, $false, $false, $false, $true, 1, $false, "98-4=94", 2) | Out-Null
$d.Content.Find.Execute("41+19=60", $true, $false, $false, $false, $false, $true, 1, $false, "95-21=74", 2) | Out-Null
$d.Content.Find.Execute("97-80=17", $true, $false, $false, $false, $false, $true, 1, $false, "74+0=74", 2) | Out-Null
$d.Content.Find.Execute("27-13=14", $true, $false, $false, $false, $false, $true, 1, $false, "12+31=43", 2) | Out-Null
$d.Content.Find.Execute("50-11=39", $true, $false, $false, $false, $false, $true, 1, $false, "59-28=31", 2) | Out-Null
$d.Content.Find.Execute("87-39=48", $true, $false, $false, $false, $false, $true, 1, $false, "60+24=84", 2) | Out-Null
$d.Content.Find.Execute("58-5=53", $true, $false, $false, $false, $false, $true, 1, $false, "34+12=46", 2) | Out-Null
$d.Content.Find.Execute("28+25=53", $true, $false, $false, $false, $false, $true, 1, $false, "30-3=27", 2) | Out-Null
$d.Content.Find.Execute("98-77=21", $true, $false, $false, $false, $false, $true, 1, $false, "77-61=16", 2) | Out-Null
$d.Content.Find.Execute("31-17=14", $true, $false, $false, $false, $false, $true, 1, $false, "82-75=7", 2) | Out-Null
$d.Content.Find.Execute("88-83=5", $true, $false, $false, $false, $false, $true, 1, $false, "81+11=92", 2) | Out-Null
$d.Content.Find.Execute("61-20=41", $true, $false, $false, $false, $false, $true, 1, $false, "70+2=72", 2) | Out-Null
$d.Content.Find.Execute("6+81=87", $true, $false, $false, $false, $false, $true, 1, $false, "77-19=58", 2) | Out-Null
$d.Content.Find.Execute("59-32=27", $true, $false, $false, $false, $false, $true, 1, $false, "1+20=21", 2) | Out-Null
$d.Content.Find.Execute("28+41=69", $true, $false, $false, $false, $false, $true, 1, $false, "32+63=95", 2) | Out-Null
$d.Content.Find.Execute("55+16=71", $true, $false, $false, $false, $false, $true, 1, $false, "58-57=1", 2) | Out-Null
$d.Content.Find.Execute("47-40=7", $true, $false, $false, $false, $false, $true, 1, $false, "93-79=14", 2) | Out-Null
$d.Content.Find.Execute("4+55=59", $true, $false, $false, $false, $false, $true, 1, $false, "24-13=11", 2) | Out-Null
$d.Content.Find.Execute("16+48=64", $true, $false, $false, $false, $false, $true, 1, $false, "21+72=93", 2) | Out-Null
$d.Content.Find.Execute("13-5=8", $true, $false, $false, $false, $false, $true, 1, $false, "7+20=27", 2) | Out-Null
$d.Content.Find.Execute("91-75=16", $true, $false, $false, $false, $false, $true, 1, $false, "7+3=10", 2) | Out-Null
$d.Content.Find.Execute("74-15=59", $true, $false, $false, $false, $false, $true, 1, $false, "91+2=93", 2) | Out-Null
$d.Content.Find.Execute("4+51=55", $true, $false, $false, $false, $false, $true, 1, $false, "56-41=15", 2) | Out-Null
$d.Content.Find.Execute("82+14=96", $true, $false, $false, $false, $false, $true, 1, $false, "2+2=4", 2) | Out-Null
$d.Content.Find.Execute("41+12=53", $true, $false, $false, $false, $false, $true, 1, $false, "45-3=42", 2) | Out-Null
$d.Content.Find.Execute("9+32=41", $true, $false, $false, $false, $false, $true, 1, $false, "44-36=8", 2) | Out-Null
$d.Content.Find.Execute("3+76=79", $true, $false, $false, $false, $false, $true, 1, $false, "81-3=78", 2) | Out-Null
$d.Content.Find.Execute("49-18=31", $true, $false, $false, $false, $false, $true, 1, $false, "71-38=33", 2) | Out-Null
$d.Content.Find.Execute("14+9=23", $true, $false, $false, $false, $false, $true, 1, $false, "74-10=64", 2) | Out-Null
$d.Content.Find.Execute("44-11=33", $true, $false, $false, $false, $false, $true, 1, $false, "2+61=63", 2) | Out-Null
$d.Content.Find.Execute("11+44=55", $true, $false, $false, $false, $false, $true, 1, $false, "7+2=9", 2) | Out-Null
$d.Content.Find.Execute("46+7=53", $true, $false, $false, $false, $false, $true, 1, $false, "15+45=60", 2) | Out-Null
$d.Content.Find.Execute("54-44=10", $true, $false, $false, $false, $false, $true, 1, $false, "61-58=3", 2) | Out-Null
$d.Content.Find.Execute("61+8=69", $true, $false, $false, $false, $false, $true, 1, $false, "87-5=82", 2) | Out-Null
$d.Content.Find.Execute("52-47=5", $true, $false, $false, $false, $false, $true, 1, $false, "95-91=4", 2) | Out-Null
